$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.299.04'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.13%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.032.81'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.35%  '

# Row 4
$ws.Range('E4').Value = '  -0.16%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.15'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.37%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.12'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.09%  '

# Row 7
$ws.Range('E7').Value = '  -0.06%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.028.92'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.20%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.522'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.09%  '

# Row 10
$ws.Range('E10').Value = '  +4.47%  '

# Row 11
$ws.Range('E11').Value = '  -1.34%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.485'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +7.13%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000247'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.80%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.75'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +7.26%  '

# Row 15
$ws.Range('E15').Value = '  +0.03%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.264.99'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.08%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.532.39'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.10%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.20'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +5.01%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.65'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +20.90%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.029.59'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.09%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '465.43'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.05%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.710'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.28%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.40'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.63%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.06'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.11%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.74'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.23%  '

# Row 26
$ws.Range('E26').Value = '  -0.18%  '

# Row 27
$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.06%  '

# Row 28
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.01'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.29%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.16'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.49%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.44'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.56%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.64'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.43%  '

# Row 32
$ws.Range('B32').Value = 'PEPE'
$ws.Range('C32').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0000101'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.00%  '

# Row 33
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.117'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.72%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '28.15'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.15%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.03%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.990'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.14%  '

# Row 37
$ws.Range('E37').Value = '  +1.38%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '48.32'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +10.29%  '

# Row 39
$ws.Range('B39').Value = 'TheGraph'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.319'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.12%  '

# Row 40
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '49.59'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.19%  '

# Row 41
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.04'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.54%  '

# Row 42
$ws.Range('E42').Value = '  -0.51%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.66'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.17%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.82'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.19%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0359'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.62%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '379.16'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.56%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.722.26'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.09%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '133.90'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.33%  '

# Row 49
$ws.Range('E49').Value = '  +0.02%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '24.47'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.07%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.22'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.64%  '
